# "corrected the typing errors"
# 1. "Request time" -> "Request Line"
# 2. The "_GoBack" bookmark (Word's auto-maintained "last edit location"
#    marker) follows the edit: it moves from after the old
#    "505 Http Variation not supported" paragraph to right after the
#    newly-edited "Request Line" text.

$d = $word.ActiveDocument

# --- 1. Fix the typo -------------------------------------------------
$d.Content.Find.Execute("Request time", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Request Line", 2) | Out-Null

# --- 2. Relocate the _GoBack bookmark to the edit point ---------------
# Locate the end of the freshly-edited text.
$editRange = $d.Content
$editRange.Find.Execute("Request Line") | Out-Null
$editEnd = $editRange.End

# Temporarily insert two placeholder characters right at the edit point
# so the collapsed bookmark range we need doesn't sit exactly on the
# paragraph-mark boundary (inserting, bookmarking between them, then
# deleting both keeps the bookmark anchored at that exact spot).
$placeholder = $d.Range($editEnd, $editEnd)
$placeholder.InsertAfter("XY")

$bookmarkRange = $d.Range($editEnd + 1, $editEnd + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

$d.Range($editEnd + 1, $editEnd + 2).Delete() | Out-Null
$d.Range($editEnd, $editEnd + 1).Delete() | Out-Null
